$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

# Set ExisUnits (column E) to 0 for rows 7 and 10
$ws.Range("E7").Value = 0
$ws.Range("E10").Value = 0

# Set MaxInvest (column S) to 100 for rows 7-11
$ws.Range("S7").Value = 100
$ws.Range("S8").Value = 100
$ws.Range("S9").Value = 100
$ws.Range("S10").Value = 100
$ws.Range("S11").Value = 100
